# Refresh the crypto price/volume table (Price = column D, Volume(1h) = column E)
# with newly scraped figures, matching the GitHub Actions data-update commit.
# Two coin pairs also swapped rank order (WrappedEther/ShibaInu at rows
# 17-18, and Bittensor/dogwifhat at rows 47-48), so their Coin/Link/Price/
# Volume cells move together as whole-row updates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '64.861.88' },
    @{ Cell = 'E2'; Value = '  -2.05%  ' },
    @{ Cell = 'D3'; Value = '3.231.19' },
    @{ Cell = 'E3'; Value = '  -1.45%  ' },
    @{ Cell = 'D5'; Value = '578.43' },
    @{ Cell = 'E5'; Value = '  +0.21%  ' },
    @{ Cell = 'D6'; Value = '173.46' },
    @{ Cell = 'E6'; Value = '  -3.27%  ' },
    @{ Cell = 'D7'; Value = '0.629' },
    @{ Cell = 'E7'; Value = '  +0.72%  ' },
    @{ Cell = 'D9'; Value = '3.227.32' },
    @{ Cell = 'E9'; Value = '  -1.43%  ' },
    @{ Cell = 'E10'; Value = '  -2.80%  ' },
    @{ Cell = 'D11'; Value = '6.78' },
    @{ Cell = 'E11'; Value = '  +1.01%  ' },
    @{ Cell = 'E12'; Value = '  -3.01%  ' },
    @{ Cell = 'D13'; Value = '3.790.00' },
    @{ Cell = 'E13'; Value = '  -1.44%  ' },
    @{ Cell = 'E14'; Value = '  -3.23%  ' },
    @{ Cell = 'D15'; Value = '64.940.73' },
    @{ Cell = 'E15'; Value = '  -2.01%  ' },
    @{ Cell = 'D16'; Value = '25.58' },
    @{ Cell = 'E16'; Value = '  -2.96%  ' },
    @{ Cell = 'B17'; Value = 'WrappedEther' },
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth' },
    @{ Cell = 'D17'; Value = '3.231.14' },
    @{ Cell = 'E17'; Value = '  -4.54%  ' },
    @{ Cell = 'B18'; Value = 'ShibaInu' },
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib' },
    @{ Cell = 'D18'; Value = '0.0000159' },
    @{ Cell = 'E18'; Value = '  -2.16%  ' },
    @{ Cell = 'D19'; Value = '415.35' },
    @{ Cell = 'E19'; Value = '  -4.45%  ' },
    @{ Cell = 'D20'; Value = '5.38' },
    @{ Cell = 'E20'; Value = '  -2.52%  ' },
    @{ Cell = 'D21'; Value = '12.84' },
    @{ Cell = 'E21'; Value = '  -1.79%  ' },
    @{ Cell = 'D22'; Value = '7.20' },
    @{ Cell = 'E22'; Value = '  -2.34%  ' },
    @{ Cell = 'D23'; Value = '1.00' },
    @{ Cell = 'E23'; Value = '  +0.20%  ' },
    @{ Cell = 'D24'; Value = '70.34' },
    @{ Cell = 'E24'; Value = '  -1.88%  ' },
    @{ Cell = 'D25'; Value = '5.65' },
    @{ Cell = 'E25'; Value = '  -0.42%  ' },
    @{ Cell = 'D26'; Value = '0.204' },
    @{ Cell = 'E26'; Value = '  +3.56%  ' },
    @{ Cell = 'D27'; Value = '0.494' },
    @{ Cell = 'E27'; Value = '  -1.87%  ' },
    @{ Cell = 'D28'; Value = '0.0000111' },
    @{ Cell = 'E28'; Value = '  -1.91%  ' },
    @{ Cell = 'D29'; Value = '9.10' },
    @{ Cell = 'E29'; Value = '  +2.83%  ' },
    @{ Cell = 'E30'; Value = '  -0.05%  ' },
    @{ Cell = 'E31'; Value = '  -2.91%  ' },
    @{ Cell = 'D32'; Value = '21.77' },
    @{ Cell = 'E32'; Value = '  -1.92%  ' },
    @{ Cell = 'E33'; Value = '  +0.06%  ' },
    @{ Cell = 'D34'; Value = '5.01' },
    @{ Cell = 'E34'; Value = '  -2.85%  ' },
    @{ Cell = 'D35'; Value = '6.40' },
    @{ Cell = 'E35'; Value = '  -2.43%  ' },
    @{ Cell = 'D36'; Value = '1.16' },
    @{ Cell = 'E36'; Value = '  -1.93%  ' },
    @{ Cell = 'D37'; Value = '156.57' },
    @{ Cell = 'E37'; Value = '  -0.65%  ' },
    @{ Cell = 'E38'; Value = '  -1.68%  ' },
    @{ Cell = 'D39'; Value = '2.823.39' },
    @{ Cell = 'E39'; Value = '  +2.41%  ' },
    @{ Cell = 'D40'; Value = '1.75' },
    @{ Cell = 'E40'; Value = '  -2.11%  ' },
    @{ Cell = 'D41'; Value = '25.43' },
    @{ Cell = 'E41'; Value = '  -4.10%  ' },
    @{ Cell = 'D42'; Value = '4.22' },
    @{ Cell = 'E42'; Value = '  -1.72%  ' },
    @{ Cell = 'D43'; Value = '0.727' },
    @{ Cell = 'E43'; Value = '  -5.83%  ' },
    @{ Cell = 'E44'; Value = '  -2.25%  ' },
    @{ Cell = 'D45'; Value = '5.75' },
    @{ Cell = 'E45'; Value = '  -4.58%  ' },
    @{ Cell = 'D46'; Value = '0.0628' },
    @{ Cell = 'E46'; Value = '  -4.13%  ' },
    @{ Cell = 'B47'; Value = 'Bittensor' },
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao' },
    @{ Cell = 'D47'; Value = '306.12' },
    @{ Cell = 'E47'; Value = '  -5.16%  ' },
    @{ Cell = 'B48'; Value = 'dogwifhat' },
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif' },
    @{ Cell = 'D48'; Value = '2.18' },
    @{ Cell = 'E48'; Value = '  -5.11%  ' },
    @{ Cell = 'D49'; Value = '22.08' },
    @{ Cell = 'E49'; Value = '  -4.70%  ' },
    @{ Cell = 'D50'; Value = '0.0263' },
    @{ Cell = 'E50'; Value = '  -0.87%  ' },
    @{ Cell = 'E51'; Value = '  -0.54%  ' }
)

# Matches plain numeric-looking text such as "578.43", "1.00", "0.0000159",
# "7.20" or "+12" (optionally signed, single decimal point).
$numericPattern = '^[+-]?[0-9]*\.?[0-9]+$'

foreach ($u in $updates) {
    $value = $u.Value
    if ($value -match $numericPattern) {
        # The sheet stores these "Price" figures as plain text (note
        # trailing zeros like "1.00" / "7.20" that a real Number would
        # lose). Assigning the bare string would let Excel auto-convert it
        # to a Number, so force text entry with a leading apostrophe -
        # exactly like a user typing '7.20 into a cell.
        $ws.Range($u.Cell).Value = "'" + $value
    } else {
        $ws.Range($u.Cell).Value = $value
    }
}
